$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1517.4445
$ws.Range("I41").Value = 962.2105
$ws.Range("K41").Value = 962.2105
$ws.Range("M41").Value = -522.2105
$ws.Range("H88").Value = 5165
$ws.Range("J88").Value = 6499
$ws.Range("L88").Value = 6499
$ws.Range("N88").Value = -7311
$ws.Range("H91").Value = 5165
$ws.Range("J91").Value = 6499
$ws.Range("L91").Value = 6499
$ws.Range("N91").Value = -9307
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = ""
$ws.Range("H132").Value = 1678.3939
$ws.Range("I132").Value = 1273.0526
$ws.Range("J132").Value = 2228.5
$ws.Range("K132").Value = 3819.1578
$ws.Range("L132").Value = 6685.5
$ws.Range("M132").Value = -1289.1578
$ws.Range("N132").Value = -11745.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6294793.5
$ws.Range("I74").Value = 9264338
$ws.Range("K74").Value = 9264338
$ws.Range("M74").Value = -9263464
$ws.Range("H77").Value = 6294793.5
$ws.Range("I77").Value = 9264338
$ws.Range("K77").Value = 46321690
$ws.Range("M77").Value = -46317322
$ws.Range("H122").Value = 2008.1777
$ws.Range("I122").Value = 1288.25
$ws.Range("K122").Value = 3864.75
$ws.Range("M122").Value = -1414.75
$ws.Range("H134").Value = 117000
$ws.Range("J134").Value = 117000
$ws.Range("L134").Value = 117000
$ws.Range("N134").Value = -127140
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 103999
$ws.Range("J53").Value = 103999
$ws.Range("L53").Value = 103999
$ws.Range("N53").Value = -105147
$ws.Range("H97").Value = 4233.5
$ws.Range("I97").Value = 4233.5
$ws.Range("K97").Value = 4233.5
$ws.Range("M97").Value = -3242.5
$ws.Range("H134").Value = 3236.7778
$ws.Range("I134").Value = 1873.25
$ws.Range("J134").Value = 4327.6
$ws.Range("K134").Value = 5619.75
$ws.Range("L134").Value = 12982.8
$ws.Range("M134").Value = -3084.75
$ws.Range("N134").Value = -18052.8
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25906.604
$ws.Range("I31").Value = 2267.6875
$ws.Range("J31").Value = 73184.44
$ws.Range("K31").Value = 2267.6875
$ws.Range("L31").Value = 73184.44
$ws.Range("M31").Value = -1972.6875
$ws.Range("N31").Value = -73774.44
$ws.Range("H34").Value = 25906.604
$ws.Range("I34").Value = 2267.6875
$ws.Range("J34").Value = 73184.44
$ws.Range("K34").Value = 2267.6875
$ws.Range("L34").Value = 73184.44
$ws.Range("M34").Value = -2065.6875
$ws.Range("N34").Value = -73588.44
$ws.Range("H132").Value = 3943.0833
$ws.Range("I132").Value = 3585.1785
$ws.Range("J132").Value = 4444.15
$ws.Range("K132").Value = 10755.5355
$ws.Range("L132").Value = 13332.45
$ws.Range("M132").Value = -8225.5355
$ws.Range("N132").Value = -18392.45
$ws.Range("H134").Value = 5939.35
$ws.Range("I134").Value = 3693.0952
$ws.Range("J134").Value = 8422.053
$ws.Range("K134").Value = 11079.2856
$ws.Range("L134").Value = 25266.159
$ws.Range("M134").Value = -8544.285600000001
$ws.Range("N134").Value = -30336.159
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 685.5714
$ws.Range("I41").Value = 299.5
$ws.Range("J41").Value = 840
$ws.Range("K41").Value = 898.5
$ws.Range("L41").Value = 2520
$ws.Range("M41").Value = -560.5
$ws.Range("N41").Value = -3196
$ws.Range("H80").Value = 7334
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").Value = ""
$ws.Range("H83").Value = 7334
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").Value = ""
$ws.Range("H92").Value = 1981.1428
$ws.Range("I92").Value = 250
$ws.Range("J92").Value = 2269.6667
$ws.Range("K92").Value = 750
$ws.Range("L92").Value = 6809.000100000001
$ws.Range("M92").Value = 498
$ws.Range("N92").Value = -9305.000100000001
$ws.Range("H97").Value = 1049.5
$ws.Range("I97").Value = 399
$ws.Range("K97").Value = 1197
$ws.Range("M97").Value = -701
$ws.Range("H107").Value = 1624.75
$ws.Range("I107").Value = 2113.375
$ws.Range("K107").Value = 6340.125
$ws.Range("M107").Value = -4420.125
$ws.Range("H109").Value = 1795
$ws.Range("J109").Value = 1700
$ws.Range("L109").Value = 5100
$ws.Range("N109").Value = -7180
$ws.Range("H122").Value = 1546.3043
$ws.Range("I122").Value = 243.83333
$ws.Range("K122").Value = 2194.49997
$ws.Range("M122").Value = 255.5000300000002
$ws.Range("H129").Value = 5954557.5
$ws.Range("I129").Value = 851.6
$ws.Range("K129").Value = 2554.8
$ws.Range("M129").Value = 2445.2
$ws.Range("H136").Value = 1113.25
$ws.Range("I136").Value = 1078.091
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 3234.273
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = 1865.727
$ws.Range("N136").Value = -14700
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3736446.5
$ws.Range("I11").Value = 3954531.8
$ws.Range("J11").Value = 29000
$ws.Range("K11").Value = 3954531.8
$ws.Range("L11").Value = 29000
$ws.Range("M11").Value = -3954392.8
$ws.Range("N11").Value = -29278
$ws.Range("H80").Value = 7240
$ws.Range("I80").Value = 5483.3335
$ws.Range("J80").Value = 9875
$ws.Range("K80").Value = 5483.3335
$ws.Range("L80").Value = 9875
$ws.Range("M80").Value = -4485.3335
$ws.Range("N80").Value = -11871
$ws.Range("H83").Value = 7240
$ws.Range("I83").Value = 5483.3335
$ws.Range("J83").Value = 9875
$ws.Range("K83").Value = 27416.6675
$ws.Range("L83").Value = 49375
$ws.Range("M83").Value = -22424.6675
$ws.Range("N83").Value = -59359
$ws.Range("H113").Value = 5572.4546
$ws.Range("I113").Value = 3999.5
$ws.Range("K113").Value = 3999.5
$ws.Range("M113").Value = -1829.5
$ws.Range("H138").Value = 75428.5
$ws.Range("J138").Value = 75428.5
$ws.Range("L138").Value = 75428.5
$ws.Range("N138").Value = -85708.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3132.7585
$ws.Range("I100").Value = 2093.84
$ws.Range("K100").Value = 2093.84
$ws.Range("M100").Value = -1552.84
$ws.Range("H114").Value = 52631
$ws.Range("J114").Value = 52631
$ws.Range("L114").Value = 52631
$ws.Range("N114").Value = -61309
$ws.Range("H132").Value = 6324.909
$ws.Range("I132").Value = 2919.8
$ws.Range("J132").Value = 9162.5
$ws.Range("K132").Value = 8759.400000000001
$ws.Range("L132").Value = 27487.5
$ws.Range("M132").Value = -6229.400000000001
$ws.Range("N132").Value = -32547.5
$ws.Range("H136").Value = 4849.6924
$ws.Range("I136").Value = 1978.0526
$ws.Range("K136").Value = 5934.1578
$ws.Range("M136").Value = -3384.1578
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 2220
$ws.Range("I7").Value = 2220
$ws.Range("K7").Value = 2220
$ws.Range("M7").Value = -2107
$ws.Range("H70").Value = 12000
$ws.Range("J70").Value = 12000
$ws.Range("L70").Value = 12000
$ws.Range("N70").Value = -12630
$ws.Range("H73").Value = 12000
$ws.Range("J73").Value = 12000
$ws.Range("L73").Value = 12000
$ws.Range("N73").Value = -14184
$ws.Range("H122").Value = 3161.3447
$ws.Range("I122").Value = 2469.4075
$ws.Range("K122").Value = 7408.2225
$ws.Range("M122").Value = -4958.2225
$ws.Range("H126").Value = 4042.7058
$ws.Range("I126").Value = 3701.8572
$ws.Range("K126").Value = 11105.5716
$ws.Range("M126").Value = -8635.571599999999
